# Generate Report for Handoff
# Marks the "aeb33cbb-bc4d-40cb-ba3e-28416ed52bfb.md" file as ready for handoff
# on both the Overview sheet and the per-language (zh-cn / de-de) detail sheets,
# updating the associated handoff timestamps.

$wb = $excel.ActiveWorkbook

$readyStatus = "Ready for handoff"

# --- Overview sheet: row 3 corresponds to aeb33cbb-bc4d-40cb-ba3e-28416ed52bfb.md ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B3").Value = $readyStatus
$overview.Range("C3").Value = $readyStatus
$overview.Range("D3").Value = "2016-03-23 06:54:52"

# --- zh-cn sheet: row 3 corresponds to aeb33cbb-bc4d-40cb-ba3e-28416ed52bfb.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $readyStatus
$zhcn.Range("E3").Value = "2016-03-23 06:54:44"

# --- de-de sheet: row 3 corresponds to aeb33cbb-bc4d-40cb-ba3e-28416ed52bfb.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $readyStatus
$dede.Range("E3").Value = "2016-03-23 06:54:52"
